$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Actual output" column (F) for rows 4 and 5
$ws.Range("F4").Value = "Its gets displayed the Supplier dropdown list"
$ws.Range("F5").Value = "Once click the Generate report button its gets displayed the details of Excel sheet"

# Update the visible selection / scroll position in the sheet view
$ws.Range("F5").Select()
$excel.ActiveWindow.ScrollColumn = 4
